# agriculture_SDMX_Model.xlsx modifications
# Splits the single "AGRICULTURE_ACTIVITY" dimension/codelist sheet into three
# separate yes/no indicator sheets: AGRICULTURE_VEGETABLE, AGRICULTURE_TUBER,
# AGRICULTURE_FRUIT - and updates the DSD sheet + codelist ids accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) DSD sheet: update dimension rows, insert the 2 extra agriculture rows
# ---------------------------------------------------------------------------
$dsd = $wb.Worksheets.Item("DSD")

# REF_AREA codelist id
$dsd.Cells.Item(4, 6).Value = "CL_COM_GEO_PICT_L123"

# INDICATOR codelist id + CDCL? flag
$dsd.Cells.Item(5, 6).Value = "CL_HH_PRIMARY_ACTIVITY_INDICATORS"
$dsd.Cells.Item(5, 7).Value = "N"

# Row 8 (was AGRICULTURE_ACTIVITY) becomes AGRICULTURE_VEGETABLE
$dsd.Cells.Item(8, 1).Value = "AGRICULTURE_VEGETABLE"
$dsd.Cells.Item(8, 2).Value = "Vegetable farming"
$dsd.Cells.Item(8, 3).Value = "N"
$dsd.Cells.Item(8, 4).Value = "Dimension"
$dsd.Cells.Item(8, 5).Value = "Coded"
$dsd.Cells.Item(8, 6).Value = "CL_COM_YESNO"
$dsd.Cells.Item(8, 7).Value = "Y"

# Insert two new rows for AGRICULTURE_TUBER and AGRICULTURE_FRUIT right after
$dsd.Range("A9:A10").EntireRow.Insert()

$dsd.Cells.Item(9, 1).Value = "AGRICULTURE_TUBER"
$dsd.Cells.Item(9, 2).Value = "Tuber farming"
$dsd.Cells.Item(9, 3).Value = "Y"
$dsd.Cells.Item(9, 4).Value = "Dimension"
$dsd.Cells.Item(9, 5).Value = "Coded"
$dsd.Cells.Item(9, 6).Value = "CL_COM_YESNO"
$dsd.Cells.Item(9, 7).Value = "Y"

$dsd.Cells.Item(10, 1).Value = "AGRICULTURE_FRUIT"
$dsd.Cells.Item(10, 2).Value = "Fruit farming"
$dsd.Cells.Item(10, 3).Value = "Y"
$dsd.Cells.Item(10, 4).Value = "Dimension"
$dsd.Cells.Item(10, 5).Value = "Coded"
$dsd.Cells.Item(10, 6).Value = "CL_COM_YESNO"
$dsd.Cells.Item(10, 7).Value = "Y"

# Column width tweaks (col A widened, col F widened)
$dsd.Columns.Item(1).ColumnWidth = 23.5
$dsd.Columns.Item(6).ColumnWidth = 36.333333333333336

$dsd.Range("F25").Select()

# ---------------------------------------------------------------------------
# 2) Agriculture codelist sheets: rename AGRICULTURE_ACTIVITY ->
#    AGRICULTURE_FRUIT, update its codelist to the shared YES/NO list, then
#    duplicate it (with formatting) to build the VEGETABLE and TUBER sheets.
# ---------------------------------------------------------------------------
$wsFruit = $wb.Worksheets.Item("AGRICULTURE_ACTIVITY")
$wsFruit.Name = "AGRICULTURE_FRUIT"

# Drop the FR/Fruit, OT/Other rows - keep header + first two data rows + Total
$wsFruit.Range("A4:A5").EntireRow.Delete()

# Relabel the remaining VG/TB data rows as the shared YES/NO codes
$wsFruit.Cells.Item(2, 1).Value = "YES"
$wsFruit.Cells.Item(2, 2).Value = "Yes"
$wsFruit.Cells.Item(3, 1).Value = "NO"
$wsFruit.Cells.Item(3, 2).Value = "No"

# Duplicate (keeps styles/col widths) twice, placed immediately before FRUIT
$wsFruit.Copy($wsFruit)
$wsFruit2 = $wb.Worksheets.Item("AGRICULTURE_FRUIT")
$wsFruit2.Copy($wsFruit2)

# After the two copies, order is: DSD, copy2, copy1, AGRICULTURE_FRUIT, Indicator
$wb.Worksheets.Item(2).Name = "AGRICULTURE_VEGETABLE"
$wb.Worksheets.Item(3).Name = "AGRICULTURE_TUBER"

$wb.Worksheets.Item("AGRICULTURE_VEGETABLE").Range("F25").Select()
$wb.Worksheets.Item("AGRICULTURE_FRUIT").Range("C13").Select()

# ---------------------------------------------------------------------------
# 3) Active sheet -> AGRICULTURE_TUBER (matches saved workbook view state)
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("AGRICULTURE_TUBER").Activate()
$wb.Worksheets.Item("AGRICULTURE_TUBER").Range("F10").Select()
